# new Madigan bike hours
# Update the Riders (C) and Average (D) columns on the Ridership sheet
# with the new bike hours data for Madigan.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$ws.Range("C2").Value = 242
$ws.Range("D2").Value = 230

$ws.Range("C3").Value = 238
$ws.Range("D3").Value = 210.72

$ws.Range("C4").Value = 244
$ws.Range("D4").Value = 210.56

$ws.Range("C5").Value = 219
$ws.Range("D5").Value = 229.71

$ws.Range("C6").Value = 241
$ws.Range("D6").Value = 239.89

$ws.Range("C7").Value = 100
$ws.Range("D7").Value = 117.84

$ws.Range("C8").Value = 61
$ws.Range("D8").Value = 99.47

$wb.Save()
